$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Countries")

# --- Rename countries (shared-string text edits) ---
$ws.Range("B7").Value = "UAE"
$ws.Range("B4").Value = "Czechia"

# --- Updated "Return" answers (replace "Undecided"/ambiguous answers with "Indifferent") ---
$ws.Range("K3").Value = "Indifferent"
$ws.Range("K10").Value = "Indifferent"
$ws.Range("K11").Value = "Indifferent"

# --- Rating corrections (Safety / Hospitality / History / Nature columns) ---
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 3
$ws.Range("H3").Value = 3
$ws.Range("H5").Value = 4
$ws.Range("I6").Value = 5
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 4
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 1
$ws.Range("G11").Value = 4
$ws.Range("H13").Value = 4
$ws.Range("J14").Value = 3
$ws.Range("H15").Value = 3
$ws.Range("G16").Value = 4
$ws.Range("H17").Value = 3

# --- Conditional formatting: 3-colour scale over the rating columns ---
$rng = $ws.Range("G2:K17")
$cf = $rng.FormatConditions.AddColorScale(3)

# --- Restore the active selection on the Countries sheet ---
$ws.Range("K5").Select()
